# Insert a new weekly data row above the current row 72 (Poroto verde, Sin especificar, Primera)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 72:97 down by one to make room for the new record
$ws.Rows.Item(72).Insert()

# Populate the new row 72 with the new weekly record.
# Same market/category info as the row that used to be at 72, new date and new volume (J).
$ws.Range("A72").Value = 1
$ws.Range("B72").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C72").Value = "Arica y Parinacota"
$ws.Range("D72").Value = 45202
$ws.Range("D72").NumberFormat = $ws.Range("D73").NumberFormat
$ws.Range("E72").Value = 15
$ws.Range("F72").Value = 100112031
$ws.Range("G72").Value = "Poroto verde"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 900
$ws.Range("K72").Value = 900
$ws.Range("L72").Value = 1000
$ws.Range("M72").Value = 950
$ws.Range("N72").Value = "$/kilo"
$ws.Range("O72").Value = "Región de Arica y Parinacota"
$ws.Range("P72").Value = 950
$ws.Range("Q72").Value = 1
$ws.Range("R72").Value = "Hortaliza"
